$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Refresh the "Number" column (D) for the existing weekly rows (2-96).
#    The values were re-pulled from source and shifted slightly for most
#    metrics (a handful of farms moved between categories); the key
#    (YearWeekIso / LastDayOfWeek / Variable) stays identical.
# ---------------------------------------------------------------------------
$dChanges = @{
    2 = 11761
    3 = 11506
    7 = 11865
    8 = 11153
    12 = 11973
    13 = 10651
    14 = 1322
    16 = 1043
    17 = 12017
    18 = 10100
    19 = 1917
    21 = 1501
    22 = 12074
    23 = 9670
    24 = 2404
    26 = 1884
    27 = 12116
    28 = 9210
    29 = 2906
    31 = 2297
    32 = 12156
    33 = 8802
    34 = 3354
    36 = 2648
    37 = 12194
    38 = 8278
    39 = 3916
    41 = 3115
    42 = 12229
    43 = 7710
    44 = 4519
    46 = 3617
    47 = 12254
    48 = 7162
    49 = 5092
    51 = 4077
    52 = 12282
    53 = 6520
    54 = 5762
    56 = 4581
    57 = 12305
    58 = 5910
    59 = 6395
    61 = 5107
    62 = 12316
    63 = 5744
    64 = 6572
    66 = 5255
    67 = 12338
    68 = 5560
    69 = 6778
    71 = 5413
    72 = 12354
    73 = 5124
    74 = 7230
    76 = 5814
    77 = 12376
    78 = 4646
    79 = 7730
    81 = 6258
    82 = 12394
    83 = 4127
    84 = 8267
    86 = 6759
    87 = 12421
    88 = 3721
    89 = 8700
    91 = 7162
    92 = 12437
    93 = 3383
    94 = 9054
    96 = 7520
}

foreach ($r in $dChanges.Keys) {
    $ws.Cells.Item($r, 4).Value = $dChanges[$r]
}

# ---------------------------------------------------------------------------
# 2. Append the new week (YearWeekIso 202507 / LastDayOfWeek 2025-02-16,
#    serial 45704) as 5 new rows, one per Variable, mirroring the layout of
#    the preceding weekly blocks.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 97;  Variable = "farms_total_count";             Number = 12458 },
    @{ Row = 98;  Variable = "farms_to_examine_count";         Number = 3071 },
    @{ Row = 99;  Variable = "farms_examined_count";           Number = 9387 },
    @{ Row = 100; Variable = "farms_examined_positive_count";  Number = 1550 },
    @{ Row = 101; Variable = "farms_examined_negative_count";  Number = 7837 }
)

# Grab the date-formatted cell style from the previous week's block (column B)
# so the new date cells keep the existing built-in short-date style (s="1")
# instead of minting a brand-new number format.
$ws.Cells.Item(96, 2).Copy() | Out-Null

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = 202507
    $ws.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 2).Value = 45704
    $ws.Cells.Item($r, 3).Value = $nr.Variable
    $ws.Cells.Item($r, 4).Value = $nr.Number
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Move the active selection to C7 (matches the saved sheetView state).
# ---------------------------------------------------------------------------
$ws.Range("C7").Select() | Out-Null
